# "Add files via upload" / "rol administrador"
# Fill in the four empty "administrador" user stories (rows 5-8) and the
# remaining blanks of the "negocio - iniciar sesion" story (row 9) on the
# "Historias de Usuario" sheet of the Product Backlog.

$wb  = $excel.ActiveWorkbook
$ws  = $wb.Worksheets.Item("Historias de Usuario")

# ---------------------------------------------------------------------
# Row 5 - HU01: Como un administrador, necesito Iniciar sesion...
# ---------------------------------------------------------------------
$ws.Range("C5").Value = "Como un administrador, necesito Iniciar sesión y poder acceder a las funcionalidades`ndel mismo, con la finalidad de autentificarme en el sistema y acceder a las funcionalidades."
$ws.Range("D5").Value = "Inicio_sesion_admin"
$ws.Range("E5").Value = "Pendiente"
$ws.Range("F5").Value = 2
$ws.Range("G5").Value = "Sprint 3"
$ws.Range("G5").WrapText = $false
$ws.Range("H5").Value = "Alta"
$ws.Range("I5").Value = "Esta historia de usuario es fundamental ya que la funcionalidad de inicio de sesión es el primer punto de contacto del administrador con el sistema."
$ws.Rows.Item(5).RowHeight = 61.2

# ---------------------------------------------------------------------
# Row 6 - HU02: Como un administrador, necesito gestionar los negocios...
# ---------------------------------------------------------------------
$ws.Range("C6").Value = "Como un administrador, necesito gestionar los negocios, con la finalidad de  registrar nuevos negocios, eliminar negocios existentes, modificar la información de los negocios y visualizas dichos negocios."
$ws.Range("D6").Value = "gestion_negocios"
$ws.Range("E6").Value = "Pendiente"
$ws.Range("F6").Value = 8
$ws.Range("G6").Value = "Sprint 4"
$ws.Range("G6").WrapText = $false
$ws.Range("H6").Value = "Alta"
$ws.Range("I6").Value = "La gestión de negocios es una funcionalidad central para los administradores, permitiendo mantener actualizada la base de datos de negocios."
$ws.Rows.Item(6).RowHeight = 58.8

# ---------------------------------------------------------------------
# Row 7 - HU03: Como un administrador, necesito Gestionar ofertas...
# ---------------------------------------------------------------------
$ws.Range("C7").Value = "Como un administrador, necesito Gestionar ofertas, con la finalidad de registrar nuevas ofertas, eliminar ofertas existentes y modificar la información de las ofertas"
$ws.Range("D7").Value = "Gestion_ofertas"
$ws.Range("E7").Value = "Pendiente"
$ws.Range("F7").Value = 8
$ws.Range("G7").Value = "Sprint 5"
$ws.Range("G7").WrapText = $false
$ws.Range("H7").Value = "Alta"
$ws.Range("I7").Value = "La gestión de ofertas es crucial para mantener actualizadas las promociones y beneficios ofrecidos a los ciudadanos."
$ws.Rows.Item(7).RowHeight = 46.2

# ---------------------------------------------------------------------
# Row 8 - HU04: Como un administrador, necesito Parametrización de greencoins...
# ---------------------------------------------------------------------
$ws.Range("C8").Value = "Como un administrador, necesito Parametrización de greencoins, con la finalidad de registrar nuevas ofertas, eliminar ofertas existentes y modificar la información de las ofertas"
$ws.Range("D8").Value = "Parametros_green"
$ws.Range("E8").Value = "Pendiente"
$ws.Range("F8").Value = 6
$ws.Range("G8").Value = "Sprint 6"
$ws.Range("G8").WrapText = $false
$ws.Range("H8").Value = "Alta"
$ws.Range("I8").Value = "Debe permitir al administrador ajustar los parámetros para adaptar el sistema a nuevas políticas o condiciones del mercado."
$ws.Rows.Item(8).RowHeight = 45

# ---------------------------------------------------------------------
# Row 9 - HU05: Como negocio, necesito iniciar sesion... (C9 unchanged)
# ---------------------------------------------------------------------
$ws.Range("D9").Value = "Inicio_sesion_negocio"
$ws.Range("E9").Value = "Pendiente"
$ws.Range("F9").Value = 2
$ws.Range("G9").Value = "Sprint 7"
$ws.Range("G9").WrapText = $false
$ws.Range("H9").Value = "Alta"
$ws.Rows.Item(9).RowHeight = 20.4

# ---------------------------------------------------------------------
# Minor row-height touch-ups elsewhere on the sheet (cosmetic, carried
# over from the resave) - keep the rest of the layout intact.
# ---------------------------------------------------------------------
$ws.Rows.Item(10).RowHeight = 43.2
$ws.Rows.Item(11).RowHeight = 42.6
$ws.Rows.Item(12).RowHeight = 46.2
$ws.Rows.Item(13).RowHeight = 64.05
$ws.Rows.Item(14).RowHeight = 67.95

# ---------------------------------------------------------------------
# Restore the cursor/selection position left behind by the editing
# session (cosmetic, mirrors the saved view state).
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Instructivo")
[void]$ws2.Activate()
[void]$ws2.Range("C12").Select()

[void]$ws.Activate()
[void]$ws.Range("C9").Select()
